$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1449.75
$ws.Range("I19").Value = 1352.625
$ws.Range("J19").Value = 1514.5
$ws.Range("K19").Value = 1352.625
$ws.Range("L19").Value = 1514.5
$ws.Range("M19").Value = -1177.625
$ws.Range("N19").Value = -1864.5
$ws.Range("H127").Value = 47619756
$ws.Range("I127").Value = 76923620
$ws.Range("J127").Value = 980
$ws.Range("K127").Value = 230770860
$ws.Range("L127").Value = 2940
$ws.Range("M127").Value = -230765900
$ws.Range("N127").Value = -12860
$ws.Range("H135").Value = 777.26086
$ws.Range("I135").Value = 767.13635
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 6904.22715
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -4369.22715
$ws.Range("N135").Value = -14070
$ws.Range("H137").Value = 2570.2222
$ws.Range("I137").Value = 1731.375
$ws.Range("K137").Value = 5194.125
$ws.Range("M137").Value = -2644.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3042.1
$ws.Range("I2").Value = 4932.6665
$ws.Range("J2").Value = 2231.8572
$ws.Range("K2").Value = 4932.6665
$ws.Range("L2").Value = 2231.8572
$ws.Range("M2").Value = -4819.6665
$ws.Range("N2").Value = -2457.8572
$ws.Range("H95").Value = 28976.25
$ws.Range("J95").Value = 28976.25
$ws.Range("L95").Value = 28976.25
$ws.Range("N95").Value = -34468.25
$ws.Range("H97").Value = 1185.1765
$ws.Range("I97").Value = 1137.3334
$ws.Range("K97").Value = 1137.3334
$ws.Range("M97").Value = -641.3334
$ws.Range("H116").Value = 3042.1
$ws.Range("I116").Value = 4932.6665
$ws.Range("J116").Value = 2231.8572
$ws.Range("K116").Value = 4932.6665
$ws.Range("L116").Value = 2231.8572
$ws.Range("M116").Value = -2638.6665
$ws.Range("N116").Value = -6819.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3042.1
$ws.Range("I3").Value = 4932.6665
$ws.Range("J3").Value = 2231.8572
$ws.Range("K3").Value = 4932.6665
$ws.Range("L3").Value = 2231.8572
$ws.Range("M3").Value = -4818.6665
$ws.Range("N3").Value = -2459.8572
$ws.Range("H94").Value = 1425
$ws.Range("I94").Value = 1497.5
$ws.Range("J94").Value = 1352.5
$ws.Range("K94").Value = 1497.5
$ws.Range("L94").Value = 1352.5
$ws.Range("M94").Value = -1046.5
$ws.Range("N94").Value = -2254.5
$ws.Range("H134").Value = 5746.7646
$ws.Range("I134").Value = 5119.5
$ws.Range("J134").Value = 6642.857
$ws.Range("K134").Value = 15358.5
$ws.Range("L134").Value = 19928.571
$ws.Range("M134").Value = -12823.5
$ws.Range("N134").Value = -24998.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1150.5714
$ws.Range("I58").Value = 1052.125
$ws.Range("J58").Value = 1465.6
$ws.Range("K58").Value = 1052.125
$ws.Range("L58").Value = 1465.6
$ws.Range("M58").Value = -849.125
$ws.Range("N58").Value = -1871.6
$ws.Range("H132").Value = 1984.6041
$ws.Range("I132").Value = 1387.1714
$ws.Range("J132").Value = 3593.077
$ws.Range("K132").Value = 4161.5142
$ws.Range("L132").Value = 10779.231
$ws.Range("M132").Value = -1631.5142
$ws.Range("N132").Value = -15839.231
$ws.Range("H134").Value = 21569992
$ws.Range("I134").Value = 2858624.2
$ws.Range("J134").Value = 62501110
$ws.Range("K134").Value = 8575872.600000001
$ws.Range("L134").Value = 187503330
$ws.Range("M134").Value = -8573337.600000001
$ws.Range("N134").Value = -187508400
$ws.Range("H136").Value = 1150.5714
$ws.Range("I136").Value = 1052.125
$ws.Range("J136").Value = 1465.6
$ws.Range("K136").Value = 3156.375
$ws.Range("L136").Value = 4396.799999999999
$ws.Range("M136").Value = -606.375
$ws.Range("N136").Value = -9496.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 14375
$ws.Range("I58").Value = 10000
$ws.Range("J58").Value = 15000
$ws.Range("K58").Value = 10000
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -9723
$ws.Range("N58").Value = -15554
$ws.Range("H80").Value = 2463.8823
$ws.Range("I80").Value = 2473.3333
$ws.Range("J80").Value = 2393
$ws.Range("K80").Value = 2473.3333
$ws.Range("L80").Value = 2393
$ws.Range("M80").Value = -1475.3333
$ws.Range("N80").Value = -4389
$ws.Range("H83").Value = 2463.8823
$ws.Range("I83").Value = 2473.3333
$ws.Range("J83").Value = 2393
$ws.Range("K83").Value = 12366.6665
$ws.Range("L83").Value = 11965
$ws.Range("M83").Value = -7374.666499999999
$ws.Range("N83").Value = -21949
$ws.Range("H102").Value = 1577.3334
$ws.Range("I102").Value = 1398.5454
$ws.Range("J102").Value = 1774
$ws.Range("K102").Value = 1398.5454
$ws.Range("L102").Value = 1774
$ws.Range("M102").Value = 223.4546
$ws.Range("N102").Value = -5018
$ws.Range("H132").Value = 2344.7173
$ws.Range("I132").Value = 2239.8684
$ws.Range("J132").Value = 2842.75
$ws.Range("K132").Value = 6719.6052
$ws.Range("L132").Value = 8528.25
$ws.Range("M132").Value = -4189.6052
$ws.Range("N132").Value = -13588.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 140762.2
$ws.Range("I93").Value = 951.5
$ws.Range("J93").Value = 233969.33
$ws.Range("K93").Value = 951.5
$ws.Range("L93").Value = 233969.33
$ws.Range("M93").Value = 296.5
$ws.Range("N93").Value = -236465.33
$ws.Range("H132").Value = 3482.111
$ws.Range("I132").Value = 3614.875
$ws.Range("J132").Value = 3216.5833
$ws.Range("K132").Value = 10844.625
$ws.Range("L132").Value = 9649.749899999999
$ws.Range("M132").Value = -8314.625
$ws.Range("N132").Value = -14709.7499
$ws.Range("H136").Value = 6803901
$ws.Range("I136").Value = 9525014
$ws.Range("J136").Value = 1120
$ws.Range("K136").Value = 28575042
$ws.Range("L136").Value = 3360
$ws.Range("M136").Value = -28572492
$ws.Range("N136").Value = -8460
$ws.Range("H139").Value = 49669.6
$ws.Range("J139").Value = 49669.6
$ws.Range("L139").Value = 49669.6
$ws.Range("N139").Value = -59949.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4654281
$ws.Range("I132").Value = 6063740.5
$ws.Range("K132").Value = 18191221.5
$ws.Range("M132").Value = -18188691.5
